# This workbook contains a weekly price log for "Poroto verde" at
# "Feria Lagunitas de Puerto Montt". The edit inserts one new daily
# record as a new row 12, shifting all the existing records (previously
# rows 12-77) down by one row (to rows 13-78).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 12; this pushes existing rows
# 12..77 down to 13..78 and keeps all of their data/styles intact.
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with the new record's data.
$ws.Cells.Item(12, 1).Value2  = 4
$ws.Cells.Item(12, 2).Value2  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(12, 3).Value2  = 'Los Lagos'
$ws.Cells.Item(12, 4).Value2  = 44701
$ws.Cells.Item(12, 5).Value2  = 10
$ws.Cells.Item(12, 6).Value2  = 100112031
$ws.Cells.Item(12, 7).Value2  = 'Poroto verde'
$ws.Cells.Item(12, 8).Value2  = 'Magnum'
$ws.Cells.Item(12, 9).Value2  = 'Primera'
$ws.Cells.Item(12, 10).Value2 = 45
$ws.Cells.Item(12, 11).Value2 = 28000
$ws.Cells.Item(12, 12).Value2 = 28000
$ws.Cells.Item(12, 13).Value2 = 28000
$ws.Cells.Item(12, 14).Value2 = '$/malla 25 kilos'
$ws.Cells.Item(12, 15).Value2 = 'Perú'
$ws.Cells.Item(12, 16).Value2 = 1120
$ws.Cells.Item(12, 17).Value2 = 25
$ws.Cells.Item(12, 18).Value2 = 'Hortaliza'
